# Update automatico via Actualizar 06-07-2020 02-58-12
# Append the newest daily record (06/06/2020) to the "Condicion_Pacientes" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row - this is what Excel does when a user types a new
# row directly below an existing table, and it keeps the table's `ref`/
# `autoFilter` ranges (and the sheet `dimension`) in sync automatically.
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the previous last row down into the new row so the
# new cells pick up the same styles (date format in column A, centered
# number format in columns B:F) instead of the workbook default style.
$ws.Range("A85:F85").Copy() | Out-Null
$ws.Range("A86:F86").PasteSpecial(-4122) | Out-Null

# Fill in the new day's values.
$ws.Range("A86").Value = 43988
$ws.Range("B86").Value = 639
$ws.Range("C86").Value = 184
$ws.Range("D86").Value = 353
$ws.Range("E86").Value = 146
$ws.Range("F86").Value = 47

# Leave the selection on the last cell entered, matching a typical manual
# data-entry flow.
$ws.Range("F86").Select() | Out-Null
